$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 2604
$ws.Range("I29").Value = 2406
$ws.Range("K29").Value = 7218
$ws.Range("M29").Value = -6937
$ws.Range("H38").Value = 106.166664
$ws.Range("I38").Value = 106.166664
$ws.Range("K38").Value = 318.499992
$ws.Range("M38").Value = 53.50000799999998
$ws.Range("H121").Value = 953044
$ws.Range("I121").Value = 265
$ws.Range("J121").Value = 1042367
$ws.Range("K121").Value = 795
$ws.Range("L121").Value = 3127101
$ws.Range("M121").Value = 952
$ws.Range("N121").Value = -3130595
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 417540.84
$ws.Range("I5").Value = 500799
$ws.Range("J5").Value = 1250
$ws.Range("K5").Value = 500799
$ws.Range("L5").Value = 1250
$ws.Range("M5").Value = -500687
$ws.Range("N5").Value = -1474
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
$ws.Range("H139").Value = 51285.5
$ws.Range("J139").Value = 51285.5
$ws.Range("L139").Value = 51285.5
$ws.Range("N139").Value = -61565.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 417540.84
$ws.Range("I4").Value = 500799
$ws.Range("J4").Value = 1250
$ws.Range("K4").Value = 500799
$ws.Range("L4").Value = 1250
$ws.Range("M4").Value = -500684
$ws.Range("N4").Value = -1480
$ws.Range("H94").Value = 648.1111
$ws.Range("I94").Value = 518
$ws.Range("J94").Value = 908.3333
$ws.Range("K94").Value = 518
$ws.Range("L94").Value = 908.3333
$ws.Range("M94").Value = -67
$ws.Range("N94").Value = -1810.3333
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 215262.2
$ws.Range("J4").Value = 215262.2
$ws.Range("L4").Value = 215262.2
$ws.Range("N4").Value = -215486.2
$ws.Range("H6").Value = 3376501.8
$ws.Range("I6").Value = 20000000
$ws.Range("J6").Value = 51802
$ws.Range("K6").Value = 20000000
$ws.Range("L6").Value = 51802
$ws.Range("M6").Value = -19999887
$ws.Range("N6").Value = -52028
$ws.Range("H7").Value = 560
$ws.Range("I7").Value = 400
$ws.Range("J7").Value = 800
$ws.Range("K7").Value = 400
$ws.Range("L7").Value = 800
$ws.Range("M7").Value = -287
$ws.Range("N7").Value = -1026
$ws.Range("H43").Value = 49653
$ws.Range("J43").Value = 49653
$ws.Range("L43").Value = 49653
$ws.Range("N43").Value = -50021
$ws.Range("H58").Value = 1854.2029
$ws.Range("I58").Value = 1575.2982
$ws.Range("J58").Value = 3179
$ws.Range("K58").Value = 1575.2982
$ws.Range("L58").Value = 3179
$ws.Range("M58").Value = -1372.2982
$ws.Range("N58").Value = -3585
$ws.Range("H92").Value = 44597
$ws.Range("J92").Value = 44597
$ws.Range("L92").Value = 44597
$ws.Range("N92").Value = -49589
$ws.Range("H95").Value = 145000
$ws.Range("J95").Value = 145000
$ws.Range("L95").Value = 145000
$ws.Range("N95").Value = -150492
$ws.Range("H99").Value = 2506
$ws.Range("I99").Value = 2506
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2506
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -1008
$ws.Range("N99").ClearContents()
$ws.Range("H101").Value = 49653
$ws.Range("J101").Value = 49653
$ws.Range("L101").Value = 49653
$ws.Range("N101").Value = -56143
$ws.Range("H104").Value = 39008
$ws.Range("J104").Value = 39008
$ws.Range("L104").Value = 39008
$ws.Range("N104").Value = -44250
$ws.Range("H106").Value = 45192
$ws.Range("J106").Value = 45192
$ws.Range("L106").Value = 45192
$ws.Range("N106").Value = -47716
$ws.Range("H108").Value = 0
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("M108").ClearContents()
$ws.Range("N108").ClearContents()
$ws.Range("H110").Value = 42999
$ws.Range("J110").Value = 42999
$ws.Range("L110").Value = 42999
$ws.Range("N110").Value = -51179
$ws.Range("H111").Value = 47628.668
$ws.Range("J111").Value = 47628.668
$ws.Range("L111").Value = 47628.668
$ws.Range("N111").Value = -55808.668
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H116").Value = 44978
$ws.Range("J116").Value = 44978
$ws.Range("L116").Value = 44978
$ws.Range("N116").Value = -54156
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
$ws.Range("H118").Value = 43942
$ws.Range("J118").Value = 43942
$ws.Range("L118").Value = 43942
$ws.Range("N118").Value = -47256
$ws.Range("H126").Value = 2506
$ws.Range("I126").Value = 2506
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 7518
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -5048
$ws.Range("N126").ClearContents()
$ws.Range("H136").Value = 1854.2029
$ws.Range("I136").Value = 1575.2982
$ws.Range("J136").Value = 3179
$ws.Range("K136").Value = 4725.8946
$ws.Range("L136").Value = 9537
$ws.Range("M136").Value = -2175.8946
$ws.Range("N136").Value = -14637
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 3375.5278
$ws.Range("J113").Value = 638.1875
$ws.Range("L113").Value = 1914.5625
$ws.Range("N113").Value = -6254.5625
$ws.Range("H120").Value = 1337666.6
$ws.Range("J120").Value = 506500
$ws.Range("L120").Value = 1519500
$ws.Range("N120").Value = -1529176
$ws.Range("H131").Value = 2527.1096
$ws.Range("I131").Value = 14879.857
$ws.Range("J131").Value = 1216.9697
$ws.Range("K131").Value = 44639.571
$ws.Range("L131").Value = 3650.9091
$ws.Range("M131").Value = -39599.571
$ws.Range("N131").Value = -13730.9091
$ws.Range("H140").Value = 1865.5938
$ws.Range("I140").Value = 1677.96
$ws.Range("J140").Value = 2535.7144
$ws.Range("K140").Value = 5033.88
$ws.Range("L140").Value = 7607.1432
$ws.Range("M140").Value = 146.1199999999999
$ws.Range("N140").Value = -17967.1432
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H108").Value = 20085.5
$ws.Range("J108").Value = 20085.5
$ws.Range("L108").Value = 20085.5
$ws.Range("N108").Value = -27765.5
$ws.Range("H119").Value = 48757
$ws.Range("J119").Value = 48757
$ws.Range("L119").Value = 48757
$ws.Range("N119").Value = -58433
$ws.Range("H135").Value = 54674
$ws.Range("J135").Value = 54674
$ws.Range("L135").Value = 54674
$ws.Range("N135").Value = -64814
$ws.Range("H137").Value = 44199
$ws.Range("J137").Value = 44199
$ws.Range("L137").Value = 44199
$ws.Range("N137").Value = -54399
$ws.Range("H138").Value = 52797.4
$ws.Range("J138").Value = 52797.4
$ws.Range("L138").Value = 52797.4
$ws.Range("N138").Value = -63077.4
$ws.Range("H139").Value = 35570
$ws.Range("J139").Value = 35570
$ws.Range("L139").Value = 35570
$ws.Range("N139").Value = -45850
$ws.Range("H141").Value = 56111.6
$ws.Range("J141").Value = 56111.6
$ws.Range("L141").Value = 56111.6
$ws.Range("N141").Value = -66471.60000000001
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2018.1177
$ws.Range("I100").Value = 1650.6666
$ws.Range("J100").Value = 2900
$ws.Range("K100").Value = 1650.6666
$ws.Range("L100").Value = 2900
$ws.Range("M100").Value = -1109.6666
$ws.Range("N100").Value = -3982
$ws.Range("H130").Value = 49429
$ws.Range("J130").Value = 49429
$ws.Range("L130").Value = 49429
$ws.Range("N130").Value = -59469

Write-Output "Applied 206 cell changes"